# FN-21 merge: extend the cancellation/modification paragraph.
#
# "...időpontja több mint 12 órára van." Ha a 12 órán belül van és le
# szeretné mondani akkor telefonon, kell értesíteni az éttermet.
#
# becomes
#
# "...időpontja több mint 12 órára van, akkor tudja csak módosítani az
# időpontot. Foglalást bármikor tud a felhasználó törölni. "

$d = $word.ActiveDocument

# 1) Right after "...12 órára van" (and before the _GoBack bookmark),
#    insert the new explanation about only being able to modify the
#    reservation in that case.
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "időpontja több mint 12 órára van", $true, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Collapse(0)
    $r1.InsertAfter(", akkor tudja csak módosítani az időpontot. Foglalást bármikor tud a felhasználó törölni.")
}

# 2) Drop the old phone-cancellation sentence, leaving a single trailing
#    space where it used to start.
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    ". Ha a 12 órán belül van és le szeretné mondani akkor telefonon, kell értesíteni az éttermet.",
    $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)
